# Scheduled-runner style refresh of market-price-derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H:N)
# across the eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1176.7
$ws.Range("I32").Value = 2245
$ws.Range("J32").Value = 909.625
$ws.Range("K32").Value = 2245
$ws.Range("L32").Value = 909.625
$ws.Range("M32").Value = -1919
$ws.Range("N32").Value = -1561.625
$ws.Range("H33").Value = 121.125
$ws.Range("I33").Value = 120
$ws.Range("J33").Value = 121.8
$ws.Range("K33").Value = 120
$ws.Range("L33").Value = 121.8
$ws.Range("M33").Value = 109
$ws.Range("N33").Value = -579.8
$ws.Range("H100").Value = 1197.375
$ws.Range("I100").Value = 864.8333
$ws.Range("K100").Value = 864.8333
$ws.Range("M100").Value = -323.8333
$ws.Range("H127").Value = 1884.1305
$ws.Range("I127").Value = 1911.5834
$ws.Range("K127").Value = 5734.7502
$ws.Range("M127").Value = -774.7502000000004
$ws.Range("H132").Value = 1106.4728
$ws.Range("I132").Value = 1004.125
$ws.Range("K132").Value = 3012.375
$ws.Range("M132").Value = -482.375
$ws.Range("H135").Value = 32258492
$ws.Range("I135").Value = 458.82144
$ws.Range("K135").Value = 4129.39296
$ws.Range("M135").Value = -1594.39296
$ws.Range("H137").Value = 1163.0741
$ws.Range("I137").Value = 789.7646999999999
$ws.Range("K137").Value = 2369.2941
$ws.Range("M137").Value = 180.7058999999999
$ws.Range("H138").Value = 1462.4783
$ws.Range("I138").Value = 1255
$ws.Range("J138").Value = 4118.2
$ws.Range("K138").Value = 3765
$ws.Range("L138").Value = 12354.6
$ws.Range("M138").Value = 1375
$ws.Range("N138").Value = -22634.6
$ws.Range("H139").Value = 50765.285
$ws.Range("J139").Value = 50765.285
$ws.Range("L139").Value = 50765.285
$ws.Range("N139").Value = -61045.285
$ws.Range("H140").Value = 60367.215
$ws.Range("J140").Value = 60367.215
$ws.Range("L140").Value = 60367.215
$ws.Range("N140").Value = -70727.215
$ws.Range("H141").Value = 701332.4399999999
$ws.Range("I141").Value = 824271.75
$ws.Range("K141").Value = 2472815.25
$ws.Range("M141").Value = -2467635.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3740.3562
$ws.Range("I32").Value = 3050.7847
$ws.Range("K32").Value = 3050.7847
$ws.Range("M32").Value = -2763.7847
$ws.Range("H33").Value = 10000
$ws.Range("I33").Value = 10000
$ws.Range("K33").Value = 10000
$ws.Range("M33").Value = -9671
$ws.Range("H61").Value = 1257.6389
$ws.Range("I61").Value = 617.96875
$ws.Range("J61").Value = 6375
$ws.Range("K61").Value = 617.96875
$ws.Range("L61").Value = 6375
$ws.Range("M61").Value = -405.96875
$ws.Range("N61").Value = -6799
$ws.Range("H92").Value = 55591.332
$ws.Range("J92").Value = 55591.332
$ws.Range("L92").Value = 55591.332
$ws.Range("N92").Value = -60583.332
$ws.Range("H110").Value = 1170.6666
$ws.Range("I110").Value = 934.7
$ws.Range("J110").Value = 2350.5
$ws.Range("K110").Value = 934.7
$ws.Range("L110").Value = 2350.5
$ws.Range("M110").Value = 1110.3
$ws.Range("N110").Value = -6440.5
$ws.Range("H132").Value = 1517.3954
$ws.Range("I132").Value = 1096.1333
$ws.Range("J132").Value = 2489.5386
$ws.Range("K132").Value = 3288.3999
$ws.Range("L132").Value = 7468.6158
$ws.Range("M132").Value = -758.3998999999999
$ws.Range("N132").Value = -12528.6158
$ws.Range("H136").Value = 1257.6389
$ws.Range("I136").Value = 617.96875
$ws.Range("J136").Value = 6375
$ws.Range("K136").Value = 1853.90625
$ws.Range("L136").Value = 19125
$ws.Range("M136").Value = 696.09375
$ws.Range("N136").Value = -24225

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6708.121
$ws.Range("I134").Value = 6536.3447
$ws.Range("J134").Value = 7953.5
$ws.Range("K134").Value = 19609.0341
$ws.Range("L134").Value = 23860.5
$ws.Range("M134").Value = -17074.0341
$ws.Range("N134").Value = -28930.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1404447.9
$ws.Range("I58").Value = 2174488.8
$ws.Range("K58").Value = 2174488.8
$ws.Range("M58").Value = -2174285.8
$ws.Range("H107").Value = 431.2857
$ws.Range("I107").Value = 436
$ws.Range("K107").Value = 436
$ws.Range("M107").Value = 1484
$ws.Range("H132").Value = 1405.9246
$ws.Range("I132").Value = 838.70734
$ws.Range("J132").Value = 3343.9167
$ws.Range("K132").Value = 2516.12202
$ws.Range("L132").Value = 10031.7501
$ws.Range("M132").Value = 13.87797999999975
$ws.Range("N132").Value = -15091.7501
$ws.Range("H134").Value = 1589.25
$ws.Range("I134").Value = 1418.079
$ws.Range("J134").Value = 2239.7
$ws.Range("K134").Value = 4254.237
$ws.Range("L134").Value = 6719.099999999999
$ws.Range("M134").Value = -1719.237
$ws.Range("N134").Value = -11789.1
$ws.Range("H136").Value = 1404447.9
$ws.Range("I136").Value = 2174488.8
$ws.Range("K136").Value = 6523466.399999999
$ws.Range("M136").Value = -6520916.399999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 959.93335
$ws.Range("H98").Value = 700.2222
$ws.Range("I98").Value = 274.75
$ws.Range("J98").Value = 1040.6
$ws.Range("K98").Value = 824.25
$ws.Range("L98").Value = 3121.8
$ws.Range("M98").Value = 673.75
$ws.Range("N98").Value = -6117.799999999999
$ws.Range("H132").Value = 1090.2858
$ws.Range("I132").Value = 799
$ws.Range("K132").Value = 7191
$ws.Range("M132").Value = -4661

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 44.57143
$ws.Range("I2").Value = 10.083333
$ws.Range("J2").Value = 90.55556
$ws.Range("K2").Value = 10.083333
$ws.Range("L2").Value = 90.55556
$ws.Range("M2").Value = 102.916667
$ws.Range("N2").Value = -316.55556
$ws.Range("H113").Value = 1650.1428
$ws.Range("J113").Value = 1300
$ws.Range("L113").Value = 1300
$ws.Range("N113").Value = -5640
$ws.Range("H132").Value = 1070903.2
$ws.Range("I132").Value = 1673933.6
$ws.Range("K132").Value = 5021800.800000001
$ws.Range("M132").Value = -5019270.800000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1020.7059
$ws.Range("I93").Value = 719.63635
$ws.Range("J93").Value = 1572.6666
$ws.Range("K93").Value = 719.63635
$ws.Range("L93").Value = 1572.6666
$ws.Range("M93").Value = 528.36365
$ws.Range("N93").Value = -4068.6666
$ws.Range("H132").Value = 1840.1621
$ws.Range("I132").Value = 1162.3334
$ws.Range("K132").Value = 3487.0002
$ws.Range("M132").Value = -957.0001999999999
$ws.Range("H136").Value = 1923.415
$ws.Range("I136").Value = 1094.4048
$ws.Range("K136").Value = 3283.2144
$ws.Range("M136").Value = -733.2143999999998
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 578.2778
$ws.Range("I113").Value = 279.2857
$ws.Range("K113").Value = 837.8571000000001
$ws.Range("M113").Value = 1332.1429
$ws.Range("H126").Value = 1942.3928
$ws.Range("I126").Value = 1006.2727
$ws.Range("J126").Value = 5374.8335
$ws.Range("K126").Value = 3018.8181
$ws.Range("L126").Value = 16124.5005
$ws.Range("M126").Value = -548.8181
$ws.Range("N126").Value = -21064.5005
$ws.Range("H132").Value = 1187.9722
$ws.Range("I132").Value = 828.5397
$ws.Range("K132").Value = 2485.6191
$ws.Range("M132").Value = 44.38090000000011
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
$ws.Range("H139").Value = 70359.8
$ws.Range("J139").Value = 70359.8
$ws.Range("L139").Value = 70359.8
$ws.Range("N139").Value = -80639.8
$ws.Range("H141").Value = 79025
$ws.Range("J141").Value = 79025
$ws.Range("L141").Value = 79025
$ws.Range("N141").Value = -89385
